$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.974.56"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.673.62"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.08%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "  +0.11%  "

# Row 5: 'BNB' -> 'BNB'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "216.07"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "

# Row 6: 'XRP' -> 'XRP'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.532"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.15%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = "  +0.09%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("E8").Value = "  +2.90%  "

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E9").Value = "  +1.91%  "

# Row 10: 'Solana' -> 'Solana'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.20"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.20%  "

# Row 11: 'TRON' -> 'TRON'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0890"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.99%  "

# Row 12: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.908.70"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.99%  "

# Row 13: 'WrappedEther' -> 'WrappedEther'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.664.78"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.58%  "

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range("E14").Value = "  +1.01%  "

# Row 15: 'Polygon' -> 'Polygon'
$ws.Range("E15").Value = "  +2.01%  "

# Row 16: 'Litecoin' -> 'Litecoin'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.75"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "

# Row 17: 'WrappedBTC' -> 'WrappedBTC'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "27.003.01"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.83%  "

# Row 18: 'BitcoinCash' -> 'BitcoinCash'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "232.32"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "

# Row 19: 'ShibaInu' -> 'ShibaInu'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0735"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "

# Row 20: 'Chainlink' -> 'Chainlink'
$ws.Range("E20").Value = "  +0.13%  "

# Row 21: 'Dai' -> 'Dai'
$ws.Range("E21").Value = "  +0.10%  "

# Row 22: 'Uniswap' -> 'Uniswap'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.45"
$c.Style = "Normal"

# Row 23: 'Toncoin' -> 'Avalanche'
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.24"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24: 'Avalanche' -> 'Toncoin'
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "

# Row 25: 'Monero' -> 'Monero'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "145.63"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26: 'Cosmos' -> 'Stellar'
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.116"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.05%  "

# Row 27: 'Stellar' -> 'Cosmos'
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.14"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "

# Row 28: 'EthereumClassic' -> 'EthereumClassic'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.85"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.43%  "

# Row 29: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("E29").Value = "  +0.03%  "

# Row 30: 'Hedera' -> 'Hedera'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0497"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.76%  "

# Row 31: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("E31").Value = "  +1.33%  "

# Row 32: 'Filecoin' -> 'Filecoin'
$ws.Range("E32").Value = "  +1.98%  "

# Row 33: 'Maker' -> 'Maker'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.455.04"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.42%  "

# Row 35: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("E35").Value = "  +4.81%  "

# Row 37: 'ARBITRUM' -> 'ARBITRUM'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.903"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.20%  "

# Row 38: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E38").Value = "  -0.59%  "

# Row 39: 'VeChain' -> 'VeChain'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0168"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "

# Row 40: 'FraxShare' -> 'FraxShare'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.05"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.44%  "

# Row 41: 'PaxDollar' -> 'PaxDollar'
$ws.Range("E41").Value = "  +0.09%  "

# Row 42: 'MXToken' -> 'MXToken'
$ws.Range("E42").Value = "  +4.28%  "

# Row 43: 'WEMIXToken' -> 'WEMIXToken'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.973"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +7.35%  "

# Row 44: 'Aave' -> 'Aave'
$ws.Range("E44").Value = "  +5.04%  "

# Row 45: 'RocketPoolETH' -> 'RocketPoolETH'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.816.76"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.89%  "

# Row 46: 'TrustWalletToken' -> 'TrustWalletToken'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.781"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47: 'Quant' -> 'Quant'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "90.67"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

# Row 48: 'RenderToken' -> 'RenderToken'
$ws.Range("E48").Value = "  +1.39%  "

# Row 49: 'Algorand' -> 'BabyDogeCoin'
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0₆0103"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "

# Row 50: 'Cronos' -> 'Algorand'
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.100"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.33%  "

# Row 51: 'EnergySwap' -> 'Cronos'
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.05%  "
